$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = $origStyle
}

Set-TextValue "D2" '64.327.17'
Set-TextValue "E2" '  +2.47%  '
Set-TextValue "D3" '2.524.67'
Set-TextValue "E3" '  +2.63%  '
Set-TextValue "E4" '  -0.02%  '
Set-TextValue "D5" '579.45'
Set-TextValue "E5" '  +1.06%  '
Set-TextValue "D6" '152.49'
Set-TextValue "E6" '  +4.29%  '
Set-TextValue "E7" '  +0.04%  '
Set-TextValue "E8" '  +0.36%  '
Set-TextValue "D9" '2.529.46'
Set-TextValue "E9" '  +2.78%  '
Set-TextValue "E10" '  +0.64%  '
Set-TextValue "E12" '  -0.22%  '
Set-TextValue "D13" '0.355'
Set-TextValue "D14" '29.40'
Set-TextValue "E14" '  +1.38%  '
Set-TextValue "E15" '  +1.74%  '
Set-TextValue "D16" '2.974.37'
Set-TextValue "D17" '64.078.72'
Set-TextValue "E17" '  +2.20%  '
Set-TextValue "D18" '2.527.19'
Set-TextValue "E18" '  +2.45%  '
Set-TextValue "E19" '  -1.08%  '
Set-TextValue "D20" '10.98'
Set-TextValue "E20" '  +0.04%  '
Set-TextValue "E21" '  +3.06%  '
Set-TextValue "D22" '328.23'
Set-TextValue "E23" '  +1.31%  '
Set-TextValue "E24" '  +0.13%  '
Set-TextValue "E25" '  -1.70%  '
Set-TextValue "D26" '65.49'
Set-TextValue "E26" '  -0.18%  '
Set-TextValue "D27" '651.37'
Set-TextValue "E27" '  -0.85%  '
Set-TextValue "E28" '  +5.34%  '
Set-TextValue "D30" '1.51'
Set-TextValue "E30" '  +4.60%  '
Set-TextValue "E31" '  +0.34%  '
Set-TextValue "D32" '8.02'
Set-TextValue "E32" '  +0.45%  '
Set-TextValue "E33" '  +1.31%  '
Set-TextValue "E34" '  +2.30%  '
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  +0.03%  '
Set-TextValue "E36" '  +1.27%  '
Set-TextValue "D37" '4.82'
Set-TextValue "E37" '  +1.46%  '
Set-TextValue "E38" '  +2.81%  '
Set-TextValue "D39" '0.372'
Set-TextValue "E39" '  +1.00%  '
Set-TextValue "D40" '18.92'
Set-TextValue "E40" '  +1.19%  '
Set-TextValue "D41" '152.66'
Set-TextValue "E41" '  +0.68%  '
Set-TextValue "E42" '  +1.92%  '
Set-TextValue "E43" '  +2.85%  '
Set-TextValue "D44" '41.97'
Set-TextValue "E44" '  +0.92%  '
Set-TextValue "D45" '162.56'
Set-TextValue "E45" '  +6.20%  '
Set-TextValue "E47" '  -2.39%  '
Set-TextValue "D48" '15.44'
Set-TextValue "E48" '  +1.36%  '
Set-TextValue "D49" '3.64'
Set-TextValue "E49" '  +1.90%  '
Set-TextValue "D50" '21.27'
Set-TextValue "E50" '  +3.86%  '
Set-TextValue "D51" '0.619'
Set-TextValue "E51" '  +2.03%  '
